# Integrate new study findings by Wei et al.
# - Rename the "Uniprot Node Hypotheses" sheet to "Uniprot Node (exclude PPI node)"
#   (this also updates the _xlnm._FilterDatabase defined name that references it).
# - Update the window/view position and the active selection on that sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Rename first worksheet; the workbook-scoped _FilterDatabase defined name
# automatically follows the new sheet name.
$ws1.Name = "Uniprot Node (exclude PPI node)"

# Make sure the renamed sheet is the active / visible one, matching tabSelected="1".
$ws1.Activate() | Out-Null

# Reposition the application window (xWindow/yWindow in workbookView).
$win = $excel.ActiveWindow
$win.Left = 1920
$win.Top = 40

# Scroll the frozen pane so row 62 is the first visible row below the freeze,
# then move the active selection to F3426.
$win.ScrollRow = 62
$win.ScrollColumn = 1
$ws1.Range("F3426").Select() | Out-Null
